$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the shared string "Objekt_ID" -> "ObjektID"
$ws.Range("A2").Value2 = "ObjektID"

# 2. Give the A2 cell (the ID field) a monospace font, matching the
#    "Liberation Mono" font stack used for ID-like values.
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.Name = "Liberation Mono;Courier New;DejaVu Sans Mono;Lucida Sans Typewriter"

# 3. Move the active selection from A5 down to row 17.
$ws.Rows(17).Select()
